$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row before current row 6 (shifts rows 6-11 down to 7-12)
$ws.Rows.Item(6).Insert()
# Insert a new header row before current row 9 (which was old row 8, now shifted)
$ws.Rows.Item(9).Insert()

# Fill in the two inserted header rows with the same header values as row 1
$ws.Range("A6").Value = "t"
$ws.Range("B6").Value = "One"
$ws.Range("C6").Value = "Two"
$ws.Range("D6").Value = "Three"

$ws.Range("A9").Value = "t"
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "Two"
$ws.Range("D9").Value = "Three"

# Add a new header-like row at the end (row 14), with a trailing-space variant in C
$ws.Range("A14").Value = "t"
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "Two "
$ws.Range("D14").Value = "Three"

# Update the selection to match the final state
$ws.Range("A14:D14").Select()
